# Apply the GA scheduling update (week-based conflict handling) to every
# "GV_Tuan_*" weekly sheet. The same cell-level change is repeated
# identically on all 15 sheets.

$wb = $excel.ActiveWorkbook

$xlPasteFormats = -4122

foreach ($ws in $wb.Worksheets) {

    # ---- Row 8 ----
    # C8: Lop CL03/vi mo/R102 -> Lop CL09/vi mo(vĩ)/R104
    $ws.Range("C8").Value = "Lớp: CL09`nMôn: Kinh tế vĩ mô`nPhòng: R104`n(Lý thuyết)"

    # E8: empty -> new class CL04/vi mo/R102 (needs the highlighted fill style)
    $ws.Range("C8").Copy()
    $ws.Range("E8").PasteSpecial($xlPasteFormats)
    $ws.Range("E8").Value = "Lớp: CL04`nMôn: Kinh tế vi mô`nPhòng: R102`n(Lý thuyết)"

    # H8: Lop CL09/vĩ mo/R103 -> Lop CL04/vĩ mo/R104
    $ws.Range("H8").Value = "Lớp: CL04`nMôn: Kinh tế vĩ mô`nPhòng: R104`n(Lý thuyết)"

    # ---- Row 9 ----
    # C9: Lop CL08/vi mo/R104 -> Lop CL08/vĩ mo/R102
    $ws.Range("C9").Value = "Lớp: CL08`nMôn: Kinh tế vĩ mô`nPhòng: R102`n(Lý thuyết)"

    # H9: Lop CL04/vĩ mo/R104 -> Lop CL08/vi mo/R104
    $ws.Range("H9").Value = "Lớp: CL08`nMôn: Kinh tế vi mô`nPhòng: R104`n(Lý thuyết)"

    # ---- Row 10 ----
    # C10: Lop CL03/vĩ mo/R103 -> Lop CL03/vĩ mo/R102
    $ws.Range("C10").Value = "Lớp: CL03`nMôn: Kinh tế vĩ mô`nPhòng: R102`n(Lý thuyết)"

    # D10: empty -> new class CL09/vi mo/R104 (previously sat in E10)
    $ws.Range("C10").Copy()
    $ws.Range("D10").PasteSpecial($xlPasteFormats)
    $ws.Range("D10").Value = "Lớp: CL09`nMôn: Kinh tế vi mô`nPhòng: R104`n(Lý thuyết)"

    # E10: CL09/vi mo/R101 -> cleared (moved to D10 above)
    $ws.Range("D8").Copy()
    $ws.Range("E10").PasteSpecial($xlPasteFormats)
    $ws.Range("E10").ClearContents()

    # G10: CL08/vĩ mo/R104 -> cleared
    $ws.Range("D8").Copy()
    $ws.Range("G10").PasteSpecial($xlPasteFormats)
    $ws.Range("G10").ClearContents()

    # ---- Row 11 ----
    # D11: CL04/vi mo/R104 -> cleared
    $ws.Range("D8").Copy()
    $ws.Range("D11").PasteSpecial($xlPasteFormats)
    $ws.Range("D11").ClearContents()

    # G11: empty -> new class CL03/vi mo/R103
    $ws.Range("C10").Copy()
    $ws.Range("G11").PasteSpecial($xlPasteFormats)
    $ws.Range("G11").Value = "Lớp: CL03`nMôn: Kinh tế vi mô`nPhòng: R103`n(Lý thuyết)"
}
